$wb = $excel.ActiveWorkbook
$cablewaySheet = $wb.ActiveSheet
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cablewaySheet)
$ws.Name = "Weir BC"
$ws.Cells.Clear()

$headers1 = @("General Information","General Information","General Information","General Information","General Information","General Information","General Information")
$headers2 = @("Station ID","Category","Station Name","Latitude","Longitude","Status","Province")
for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers1[$c-1]
    $ws.Cells.Item(2, $c).Value = $headers2[$c-1]
}

$data = @(
    @("08NE087", "Weir", "DEER CREEK AT DEER PARK", "49.44847", "-118.042", "ACTIVE", "BC"),
    @("08NH115", "Weir", "SULLIVAN CREEK NEAR CANYON", "49.10427", "-116.43675", "ACTIVE", "BC"),
    @("08NJ130", "Weir", "ANDERSON CREEK NEAR NELSON", "49.50195", "-117.26128", "ACTIVE", "BC"),
    @("08HA016", "Weir", "BINGS CREEK NEAR THE MOUTH", "48.78922", "-123.72547", "ACTIVE", "BC"),
    @("08HB048", "Weir", "CARNATION CREEK AT THE MOUTH", "48.91556", "-124.99778", "ACTIVE", "BC"),
    @("08GA061", "Weir", "MACKAY CREEK AT MONTROYAL BOULEVARD", "49.35616", "-123.09986", "ACTIVE", "BC"),
    @("08MH152", "Weir", "BERTRAND CREEK AT INTERNATIONAL BOUNDARY", "49.00239", "-122.5234", "ACTIVE", "BC"),
    @("08OA003", "Weir", "PREMIER CREEK NEAR DAAJING GIIDS", "53.25861", "-132.07505", "ACTIVE", "BC"),
    @("08EE025", "Weir", "TWO MILE CREEK IN DISTRICT LOT 4834", "55.29611", "-127.62075", "ACTIVE", "BC"),
    @("08JB013", "Weir", "NORTH BEACH CREEK ABOVE ALLIN CREEK", "54.12694", "-125.92833", "ACTIVE", "BC"),
    @("08LB024", "Weir", "FISHTRAP CREEK NEAR MCLURE", "51.12203", "-120.20894", "ACTIVE", "BC"),
    @("08LE108", "Weir", "EAST CANOE CREEK ABOVE DAM", "50.69665", "-119.19616", "ACTIVE", "BC"),
    @("08LF007", "Weir", "CRISS CREEK NEAR SAVONA", "50.8839", "-120.96594", "ACTIVE", "BC"),
    @("08LG056", "Weir", "GUICHON CREEK ABOVE TUNKWA LAKE DIVERSION", "50.60783", "-120.91086", "ACTIVE", "BC"),
    @("08NM065", "Weir", "VERNON CREEK AT OUTLET OF KALAMALKA LAKE", "50.23847", "-119.26689", "ACTIVE", "BC"),
    @("08NM134", "Weir", "CAMP CREEK AT MOUTH NEAR THIRSK", "49.71234", "-120.01117", "ACTIVE", "BC"),
    @("08NM139", "Weir", "ESPERON CREEK NEAR KELOWNA", "50.07333", "-119.68778", "ACTIVE", "BC"),
    @("08NM146", "Weir", "CLARK CREEK NEAR WINFIELD", "50.04875", "-119.33597", "ACTIVE", "BC"),
    @("08NM165", "Weir", "LAMBLY CREEK ABOVE TERRACE CREEK", "49.99417", "-119.61444", "ACTIVE", "BC"),
    @("08NM171", "Weir", "VASEUX CREEK ABOVE SOLCO CREEK", "49.24928", "-119.32139", "ACTIVE", "BC"),
    @("08NM172", "Weir", "PEARSON CREEK NEAR THE MOUTH", "49.88694", "-119.06167", "ACTIVE", "BC"),
    @("08NM173", "Weir", "GREATA CREEK NEAR THE MOUTH", "49.79458", "-119.85203", "ACTIVE", "BC"),
    @("08NM240", "Weir", "TWO FORTY CREEK NEAR PENTICTON", "49.65089", "-119.40003", "ACTIVE", "BC"),
    @("08NM241", "Weir", "TWO FORTY-ONE CREEK NEAR PENTICTON", "49.65004", "-119.39385", "ACTIVE", "BC"),
    @("08NM242", "Weir", "DENNIS CREEK NEAR 1780 METRE CONTOUR", "49.62417", "-119.38143", "ACTIVE", "BC"),
    @("07SB003", "Weir", "YELLOWKNIFE RIVER AT INLET TO PROSPEROUS LAKE", "62.672", "-114.26122", "ACTIVE", "BC"),
    @("07SB012", "Weir", "DUNCAN LAKE NEAR YELLOWKNIFE", "62.80781", "-114.03336", "ACTIVE", "BC"),
)

$rowCount = $data.Count
# Latitude/Longitude columns (D, E) must stay text so numeric-looking values
# are not auto-converted to numbers by Excel.
$latLonRange = $ws.Range("D3:E" + (2 + $rowCount))
$latLonRange.NumberFormat = "@"

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws.Cells.Item($r + 3, $c + 1).Value = $data[$r][$c]
    }
}

$latLonRange.ClearFormats()

Write-Output "done"